# Insert a new data row at row 78 (pushing existing rows 78..147 down to 79..148)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 78, shifting rows 78-147 to 79-148.
$ws.Rows.Item(78).Insert()

# Populate the new row 78 with the constant columns shared by every data row
# in this sheet, plus the new record's specific values.
$ws.Cells.Item(78, 1).Value = 4                                            # A - Mercado ID
$ws.Cells.Item(78, 2).Value = "Feria Lagunitas de Puerto Montt"            # B - Mercado
$ws.Cells.Item(78, 3).Value = "Los Lagos"                                  # C - Region
$ws.Cells.Item(78, 4).Value = 44893                                        # D - Fecha
$ws.Cells.Item(78, 5).Value = 10                                           # E - Codreg
$ws.Cells.Item(78, 6).Value = 100112022                                    # F - Categoria ID
$ws.Cells.Item(78, 7).Value = "Arveja Verde"                               # G - Categoria
$ws.Cells.Item(78, 8).Value = "Sin especificar"                            # H - Variedad
$ws.Cells.Item(78, 9).Value = "Primera"                                    # I - Calidad
$ws.Cells.Item(78, 10).Value = 70                                          # J - Volumen
$ws.Cells.Item(78, 11).Value = 35000                                       # K - Precio minimo
$ws.Cells.Item(78, 12).Value = 35000                                       # L - Precio maximo
$ws.Cells.Item(78, 13).Value = 35000                                       # M - Precio promedio ponderado
$ws.Cells.Item(78, 14).Value = "$/saco 25 kilos"                           # N - Unidad de comercializacion
$ws.Cells.Item(78, 15).Value = "Región de La Araucanía"                    # O - Origen
$ws.Cells.Item(78, 16).Value = 1400                                        # P - Precio $/Kg
$ws.Cells.Item(78, 17).Value = 25                                          # Q - Kg o Unidades
$ws.Cells.Item(78, 18).Value = "Hortaliza"                                 # R - Clasificacion

# Match the date format used by the other rows in column D.
$ws.Cells.Item(78, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
